$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I19").Value = 0.03477683984818711
$ws.Range("J19").Value = 0.2704625831437422
$ws.Range("K19").Value = 0.2362491449239708
$ws.Range("L19").Value = 2.364507701123102
